$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink annotation from G2 before we shuffle data
# around (the link target cell moves from G2 to I2). Hyperlinks.Delete()
# leaves the "Hyperlink" cell style behind, so reset it back to Normal.
$ws.Range("G2").Hyperlinks.Delete()
$ws.Range("G2").Style = "Normal"

# --- Header row (row 1) ----------------------------------------------------
$ws.Range("A1").Value = "nama_lengkap"
$ws.Range("B1").Value = "alamat_mitra"
$ws.Range("C1").Value = "kode_desa"
$ws.Range("D1").Value = "kode_kecamatan"
$ws.Range("E1").Value = "kode_kabupaten"
$ws.Range("F1").Value = "kode_provinsi"
$ws.Range("G1").Value = "jenis_kelamin"
$ws.Range("H1").Value = "no_hp_mitra"
$ws.Range("I1").Value = "email_mitra"
$ws.Range("J1").Value = "tahun"

# --- Data row (row 2) -------------------------------------------------------
# A2/B2/C2/D2 stay put; the old E2 ("kelamin"=1) and F2 ("hp") values shift
# two columns right to make room for the new kode_kabupaten / kode_provinsi
# columns, and the email moves from G2 to I2.
$ws.Range("E2").Clear()
$ws.Range("F2").Clear()

$ws.Range("G2").Value = 1
$ws.Range("G2").NumberFormat = "@"

$ws.Range("H2").Value = "+62 81"
$ws.Range("I2").Value = "email@gmail.com"

# Re-create the hyperlink on its new home cell, restoring the Hyperlink
# style/number-format that the source cell had.
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:email@gmail.com") | Out-Null
$ws.Range("I2").Style = "Hyperlink"
$ws.Range("I2").NumberFormat = "@"

# --- View / selection tidy-up ------------------------------------------------
$ws.Range("K3").Select()
